$wb = $excel.ActiveWorkbook

# The first sheet (sheetId 2 / rId1) is "OperatiiExistente" -> rename to "ImplementedOperations"
$ws = $wb.Worksheets.Item("OperatiiExistente")
$ws.Name = "ImplementedOperations"

# F1 currently holds "Variants" -> change to "Summary"
$ws.Range("F1").Value = "Summary"

# Clear the helper formulas that were in E7:E9 (C+D*P2)
$ws.Range("E7").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("E9").ClearContents()

# Column H narrower, column I now hidden
$ws.Columns("H").ColumnWidth = 36.28515625
$ws.Columns("I").Hidden = $true

# Update the last active selection on the sheet
$ws.Range("F16").Select()
